$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2-10
# from 45185 (2023-09-16) to 45204 (2023-10-05)
$ws.Range("C2:C10").Value = 45204
